$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 11122.714
$ws.Range("J32").Value = 9395.875
$ws.Range("L32").Value = 9395.875
$ws.Range("N32").Value = -10047.875
$ws.Range("H98").Value = 1286.3684
$ws.Range("I98").Value = 1382.4706
$ws.Range("J98").Value = 469.5
$ws.Range("K98").Value = 1382.4706
$ws.Range("L98").Value = 469.5
$ws.Range("M98").Value = 115.5293999999999
$ws.Range("N98").Value = -3465.5
$ws.Range("H122").Value = 1286.3684
$ws.Range("I122").Value = 1382.4706
$ws.Range("J122").Value = 469.5
$ws.Range("K122").Value = 4147.4118
$ws.Range("L122").Value = 1408.5
$ws.Range("M122").Value = -1697.4118
$ws.Range("N122").Value = -6308.5
$ws.Range("H138").Value = 5297
$ws.Range("I138").Value = 5297
$ws.Range("K138").Value = 15891
$ws.Range("M138").Value = -10751
$ws.Range("H141").Value = 9782.833000000001
$ws.Range("J141").Value = 1405
$ws.Range("L141").Value = 4215
$ws.Range("N141").Value = -14575

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4907.6
$ws.Range("I45").Value = 4631
$ws.Range("K45").Value = 4631
$ws.Range("M45").Value = -4254
$ws.Range("H74").Value = 504.5
$ws.Range("I74").Value = 504.5
$ws.Range("K74").Value = 504.5
$ws.Range("M74").Value = 369.5
$ws.Range("H77").Value = 504.5
$ws.Range("I77").Value = 504.5
$ws.Range("K77").Value = 2522.5
$ws.Range("M77").Value = 1845.5
$ws.Range("H102").Value = 2040.3334
$ws.Range("I102").Value = 2040.3334
$ws.Range("K102").Value = 2040.3334
$ws.Range("M102").Value = -418.3334
$ws.Range("H140").Value = 49998
$ws.Range("J140").Value = 49998
$ws.Range("L140").Value = 49998
$ws.Range("N140").Value = -60358

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 17666.666
$ws.Range("I19").Value = 8500
$ws.Range("J19").Value = 25000
$ws.Range("K19").Value = 8500
$ws.Range("L19").Value = 25000
$ws.Range("M19").Value = -8327
$ws.Range("N19").Value = -25346
$ws.Range("H86").Value = 9314.611000000001
$ws.Range("I86").Value = 9848.571
$ws.Range("K86").Value = 9848.571
$ws.Range("M86").Value = -8725.571
$ws.Range("H88").Value = 24343
$ws.Range("J88").Value = 24343
$ws.Range("L88").Value = 24343
$ws.Range("N88").Value = -25155
$ws.Range("H89").Value = 9314.611000000001
$ws.Range("I89").Value = 9848.571
$ws.Range("K89").Value = 49242.855
$ws.Range("M89").Value = -43626.855
$ws.Range("H91").Value = 24343
$ws.Range("J91").Value = 24343
$ws.Range("L91").Value = 24343
$ws.Range("N91").Value = -27151
$ws.Range("H140").Value = 99499.5
$ws.Range("J140").Value = 99499.5
$ws.Range("L140").Value = 99499.5
$ws.Range("N140").Value = -109859.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2201.6667
$ws.Range("I16").Value = 639.4
$ws.Range("K16").Value = 639.4
$ws.Range("M16").Value = -352.4
$ws.Range("H22").Value = 4446055.5
$ws.Range("I22").Value = 1499.7142
$ws.Range("K22").Value = 1499.7142
$ws.Range("M22").Value = -1149.7142
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("H94").Value = 694
$ws.Range("J94").Value = 694
$ws.Range("L94").Value = 694
$ws.Range("N94").Value = -1596
$ws.Range("H113").Value = 2201.6667
$ws.Range("I113").Value = 639.4
$ws.Range("K113").Value = 639.4
$ws.Range("M113").Value = 1530.6
$ws.Range("N81").ClearContents()
$ws.Range("N84").ClearContents()

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2369104.2
$ws.Range("I4").Value = 902277.2
$ws.Range("K4").Value = 2706831.6
$ws.Range("M4").Value = -2706719.6
$ws.Range("H139").Value = 9258.875
$ws.Range("I139").Value = 9258.875
$ws.Range("K139").Value = 27776.625
$ws.Range("M139").Value = -22636.625
$ws.Range("H140").Value = 1003947.4
$ws.Range("I140").Value = 1003947.4
$ws.Range("K140").Value = 3011842.2
$ws.Range("M140").Value = -3006662.2

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 13079989
$ws.Range("J11").Value = 2237833.2
$ws.Range("L11").Value = 2237833.2
$ws.Range("N11").Value = -2238111.2
$ws.Range("H12").Value = 1616667
$ws.Range("I12").Value = 1500000
$ws.Range("K12").Value = 1500000
$ws.Range("M12").Value = -1499860
$ws.Range("H43").Value = 9567.733
$ws.Range("I43").Value = 724.1111
$ws.Range("K43").Value = 724.1111
$ws.Range("M43").Value = -573.1111
$ws.Range("H46").Value = 30723.75
$ws.Range("J46").Value = 30723.75
$ws.Range("L46").Value = 30723.75
$ws.Range("N46").Value = -31035.75
$ws.Range("H70").Value = 7675
$ws.Range("I70").Value = 7400.8335
$ws.Range("J70").Value = 8497.5
$ws.Range("K70").Value = 7400.8335
$ws.Range("L70").Value = 8497.5
$ws.Range("M70").Value = -7130.8335
$ws.Range("N70").Value = -9037.5
$ws.Range("H73").Value = 7675
$ws.Range("I73").Value = 7400.8335
$ws.Range("J73").Value = 8497.5
$ws.Range("K73").Value = 7400.8335
$ws.Range("L73").Value = 8497.5
$ws.Range("M73").Value = -6464.8335
$ws.Range("N73").Value = -10369.5
$ws.Range("H122").Value = 4554.6665
$ws.Range("I122").Value = 3798.8
$ws.Range("J122").Value = 5499.5
$ws.Range("K122").Value = 11396.4
$ws.Range("L122").Value = 16498.5
$ws.Range("M122").Value = -8946.400000000001
$ws.Range("N122").Value = -21398.5
$ws.Range("H134").Value = 27629.334
$ws.Range("J134").Value = 27629.334
$ws.Range("L134").Value = 82888.00199999999
$ws.Range("N134").Value = -87958.00199999999

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 344.85715
$ws.Range("I16").Value = 344.85715
$ws.Range("K16").Value = 344.85715
$ws.Range("M16").Value = -174.85715
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("H55").Value = 554.2727
$ws.Range("I55").Value = 349.66666
$ws.Range("K55").Value = 349.66666
$ws.Range("M55").Value = -176.66666
$ws.Range("N17").ClearContents()
